# Auto-generated edit script applying the market-price / profit data
# refresh captured in the commit diff ("chore: update Sheets via scheduled runner").
# Each leve row's currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns
# (H:N) are overwritten with refreshed values; one row (GSM!111) drops its
# HQ-profit cell (N111) entirely because the refreshed HQ price is 0.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 758.0526
$ws.Range("I28").Value = 787.875
$ws.Range("K28").Value = 787.875
$ws.Range("M28").Value = -302.875
$ws.Range("H98").Value = 4239.6665
$ws.Range("I98").Value = 2331.8462
$ws.Range("J98").Value = 9200
$ws.Range("K98").Value = 2331.8462
$ws.Range("L98").Value = 9200
$ws.Range("M98").Value = -833.8462
$ws.Range("N98").Value = -12196
$ws.Range("H122").Value = 4239.6665
$ws.Range("I122").Value = 2331.8462
$ws.Range("J122").Value = 9200
$ws.Range("K122").Value = 6995.5386
$ws.Range("L122").Value = 27600
$ws.Range("M122").Value = -4545.5386
$ws.Range("N122").Value = -32500
$ws.Range("H129").Value = 962.3043
$ws.Range("J129").Value = 996.907
$ws.Range("L129").Value = 2990.721
$ws.Range("N129").Value = -12990.721
$ws.Range("H132").Value = 27137072
$ws.Range("I132").Value = 27890462
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 83671386
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -83668856
$ws.Range("N132").Value = -50060
$ws.Range("H137").Value = 4341.114
$ws.Range("I137").Value = 4473.769
$ws.Range("J137").Value = 4149.5
$ws.Range("K137").Value = 13421.307
$ws.Range("L137").Value = 12448.5
$ws.Range("M137").Value = -10871.307
$ws.Range("N137").Value = -17548.5
$ws.Range("H138").Value = 4102.5
$ws.Range("I138").Value = 2287.4
$ws.Range("J138").Value = 4369.4263
$ws.Range("K138").Value = 6862.200000000001
$ws.Range("L138").Value = 13108.2789
$ws.Range("M138").Value = -1722.200000000001
$ws.Range("N138").Value = -23388.2789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14423.718
$ws.Range("I32").Value = 9367.857
$ws.Range("K32").Value = 9367.857
$ws.Range("M32").Value = -9080.857
$ws.Range("H74").Value = 5534.96
$ws.Range("I74").Value = 6791.5
$ws.Range("J74").Value = 3301.111
$ws.Range("K74").Value = 6791.5
$ws.Range("L74").Value = 3301.111
$ws.Range("M74").Value = -5917.5
$ws.Range("N74").Value = -5049.111
$ws.Range("H77").Value = 5534.96
$ws.Range("I77").Value = 6791.5
$ws.Range("J77").Value = 3301.111
$ws.Range("K77").Value = 33957.5
$ws.Range("L77").Value = 16505.555
$ws.Range("M77").Value = -29589.5
$ws.Range("N77").Value = -25241.555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 21740858
$ws.Range("I94").Value = 31251578
$ws.Range("K94").Value = 31251578
$ws.Range("M94").Value = -31251127
$ws.Range("H134").Value = 2941.98
$ws.Range("I134").Value = 1463.3889
$ws.Range("J134").Value = 6744.0713
$ws.Range("K134").Value = 4390.1667
$ws.Range("L134").Value = 20232.2139
$ws.Range("M134").Value = -1855.1667
$ws.Range("N134").Value = -25302.2139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1611
$ws.Range("I16").Value = 1518.5
$ws.Range("J16").Value = 1796
$ws.Range("K16").Value = 1518.5
$ws.Range("L16").Value = 1796
$ws.Range("M16").Value = -1231.5
$ws.Range("N16").Value = -2370
$ws.Range("H31").Value = 4822.9487
$ws.Range("I31").Value = 2450
$ws.Range("J31").Value = 5435.3228
$ws.Range("K31").Value = 2450
$ws.Range("L31").Value = 5435.3228
$ws.Range("M31").Value = -2155
$ws.Range("N31").Value = -6025.3228
$ws.Range("H34").Value = 4822.9487
$ws.Range("I34").Value = 2450
$ws.Range("J34").Value = 5435.3228
$ws.Range("K34").Value = 2450
$ws.Range("L34").Value = 5435.3228
$ws.Range("M34").Value = -2248
$ws.Range("N34").Value = -5839.3228
$ws.Range("H58").Value = 2185.879
$ws.Range("I58").Value = 1741.8474
$ws.Range("K58").Value = 1741.8474
$ws.Range("M58").Value = -1538.8474
$ws.Range("H113").Value = 1611
$ws.Range("I113").Value = 1518.5
$ws.Range("J113").Value = 1796
$ws.Range("K113").Value = 1518.5
$ws.Range("L113").Value = 1796
$ws.Range("M113").Value = 651.5
$ws.Range("N113").Value = -6136
$ws.Range("H132").Value = 3736.9429
$ws.Range("I132").Value = 3188.1667
$ws.Range("J132").Value = 4934.273
$ws.Range("K132").Value = 9564.500100000001
$ws.Range("L132").Value = 14802.819
$ws.Range("M132").Value = -7034.500100000001
$ws.Range("N132").Value = -19862.819
$ws.Range("H136").Value = 2185.879
$ws.Range("I136").Value = 1741.8474
$ws.Range("K136").Value = 5225.5422
$ws.Range("M136").Value = -2675.5422

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1968.125
$ws.Range("I5").Value = 582
$ws.Range("J5").Value = 4278.3335
$ws.Range("K5").Value = 1746
$ws.Range("L5").Value = 12835.0005
$ws.Range("M5").Value = -1634
$ws.Range("N5").Value = -13059.0005
$ws.Range("H107").Value = 1397.1765
$ws.Range("J107").Value = 2425.25
$ws.Range("L107").Value = 7275.75
$ws.Range("N107").Value = -11115.75
$ws.Range("H113").Value = 567.6739
$ws.Range("I113").Value = 575.92
$ws.Range("J113").Value = 557.8570999999999
$ws.Range("K113").Value = 1727.76
$ws.Range("L113").Value = 1673.5713
$ws.Range("M113").Value = 442.2400000000002
$ws.Range("N113").Value = -6013.5713
$ws.Range("H122").Value = 2625.9194
$ws.Range("J122").Value = 2979.0386
$ws.Range("L122").Value = 26811.3474
$ws.Range("N122").Value = -31711.3474
$ws.Range("H131").Value = 8200985.5
$ws.Range("J131").Value = 985.2826
$ws.Range("L131").Value = 2955.8478
$ws.Range("N131").Value = -13035.8478
$ws.Range("H132").Value = 1764.1538
$ws.Range("J132").Value = 2069.6155
$ws.Range("L132").Value = 18626.5395
$ws.Range("N132").Value = -23686.5395
$ws.Range("H135").Value = 1968.125
$ws.Range("I135").Value = 582
$ws.Range("J135").Value = 4278.3335
$ws.Range("K135").Value = 5238
$ws.Range("L135").Value = 38505.0015
$ws.Range("M135").Value = -2703
$ws.Range("N135").Value = -43575.0015
$ws.Range("H137").Value = 8613.954
$ws.Range("I137").Value = 3056.875
$ws.Range("J137").Value = 23432.834
$ws.Range("K137").Value = 9170.625
$ws.Range("L137").Value = 70298.50199999999
$ws.Range("M137").Value = -4070.625
$ws.Range("N137").Value = -80498.50199999999
$ws.Range("H140").Value = 21902.28
$ws.Range("I140").Value = 30797.47
$ws.Range("K140").Value = 92392.41
$ws.Range("M140").Value = -87212.41
$ws.Range("H141").Value = 7421.75
$ws.Range("I141").Value = 7013.4546
$ws.Range("K141").Value = 21040.3638
$ws.Range("M141").Value = -15860.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 2596
$ws.Range("I122").Value = 1897.1154
$ws.Range("J122").Value = 6230.2
$ws.Range("K122").Value = 5691.3462
$ws.Range("L122").Value = 18690.6
$ws.Range("M122").Value = -3241.3462
$ws.Range("N122").Value = -23590.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7278.091
$ws.Range("I40").Value = 7169.875
$ws.Range("J40").Value = 7566.6665
$ws.Range("K40").Value = 7169.875
$ws.Range("L40").Value = 7566.6665
$ws.Range("M40").Value = -7033.875
$ws.Range("N40").Value = -7838.6665
$ws.Range("H110").Value = 26622.25
$ws.Range("J110").Value = 26622.25
$ws.Range("L110").Value = 26622.25
$ws.Range("N110").Value = -34802.25
$ws.Range("H122").Value = 3669.2058
$ws.Range("I122").Value = 2514.261
$ws.Range("J122").Value = 6084.091
$ws.Range("K122").Value = 7542.782999999999
$ws.Range("L122").Value = 18252.273
$ws.Range("M122").Value = -5092.782999999999
$ws.Range("N122").Value = -23152.273
$ws.Range("H132").Value = 4753.686
$ws.Range("I132").Value = 2314.375
$ws.Range("J132").Value = 8862
$ws.Range("K132").Value = 6943.125
$ws.Range("L132").Value = 26586
$ws.Range("M132").Value = -4413.125
$ws.Range("N132").Value = -31646

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 685.7
$ws.Range("I107").Value = 482.26666
$ws.Range("J107").Value = 1296
$ws.Range("K107").Value = 1446.79998
$ws.Range("L107").Value = 3888
$ws.Range("M107").Value = 473.20002
$ws.Range("N107").Value = -7728
$ws.Range("H113").Value = 9426.272000000001
$ws.Range("I113").Value = 12773.875
$ws.Range("J113").Value = 499.33334
$ws.Range("K113").Value = 38321.625
$ws.Range("L113").Value = 1498.00002
$ws.Range("M113").Value = -36151.625
$ws.Range("N113").Value = -5838.000019999999
$ws.Range("H132").Value = 6539818.5
$ws.Range("I132").Value = 4443.3667
$ws.Range("J132").Value = 15876069
$ws.Range("K132").Value = 13330.1001
$ws.Range("L132").Value = 47628207
$ws.Range("M132").Value = -10800.1001
$ws.Range("N132").Value = -47633267
$ws.Range("H136").Value = 4826.737
$ws.Range("I136").Value = 1640.7
$ws.Range("J136").Value = 8366.777
$ws.Range("K136").Value = 4922.1
$ws.Range("L136").Value = 25100.331
$ws.Range("M136").Value = -2372.1
$ws.Range("N136").Value = -30200.331
